$wb = $excel.ActiveWorkbook

# Overview sheet: update Latest HO Xliff Generate Date for e11a4c50 row (row 3, col G)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-27 16:47:40"

# zh-cn sheet: update Correspond Handoff Datetime (H3) and Correspond Handback DateTime (K3) for e11a4c50 row
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-08-27 16:47:35"
$wsZhCn.Range("K3").Value = "2016-08-27 16:47:52"

# de-de sheet: update Correspond Handoff Datetime (H3) and Correspond Handback DateTime (K3) for e11a4c50 row
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-08-27 16:47:40"
$wsDeDe.Range("K3").Value = "2016-08-27 16:47:59"
